# "reload + ammo counter" - log time for 2024-03-14 (row 34) against the
# reload mechanic / ammo counter work, and repoint the "remaining days"
# estimate (I9) at the just-filled-in row instead of the previous blank one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LO1 (utilities): 7 minutes
$ws.Range("B34").Formula = "=(1/60)*(7)"
# LO2 (Character + NPC): 7 minutes
$ws.Range("C34").Formula = "=(1/60)*(7)"
# LO4 (world interactions): 22+22+8 minutes
$ws.Range("E34").Formula = "=(1/60)*(22+22+8)"

# "BASED ON LAST DAY" now divides by the latest logged day's total (F34)
# instead of the prior one (F33).
$ws.Range("I9").Formula = "=I3/F34"

# Leave the selection on the cell that was just edited.
$ws.Range("B35").Select()
